$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

# Fill column A ("cs") for all new rows in one shot
$ws.Range("A719:A735").Value = "cs"

# Fill the new translation key/value rows (719-735)
$ws.Range("B719").Value = "inventory.index.title"
$ws.Range("C719").Value = "Inventář"
$ws.Range("B720").Value = "inventory.home.title"
$ws.Range("C720").Value = "Inventář"
$ws.Range("B721").Value = "inventory.home.subtitle"
$ws.Range("C721").Value = "V této sekci naleznete veškeré předměty, které jste si pořídili na tržišti."
$ws.Range("B722").Value = "inventory.home.menu"
$ws.Range("C722").Value = "Inventář"
$ws.Range("B723").Value = "inventory.atomizer.menu"
$ws.Range("C723").Value = "Atomizéry"
$ws.Range("B724").Value = "inventory.mod.menu"
$ws.Range("C724").Value = "Mody"
$ws.Range("B725").Value = "inventory.cell.menu"
$ws.Range("C725").Value = "Články"
$ws.Range("B726").Value = "inventory.cotton.menu"
$ws.Range("C726").Value = "Vaty"
$ws.Range("B727").Value = "inventory.aroma.menu"
$ws.Range("C727").Value = "Aromata"
$ws.Range("B728").Value = "inventory.base.menu"
$ws.Range("C728").Value = "Báze"
$ws.Range("B729").Value = "inventory.booster.menu"
$ws.Range("C729").Value = "Boostery"
$ws.Range("B730").Value = "inventory.lab.menu"
$ws.Range("C730").Value = "Laboratoř"
$ws.Range("B731").Value = "inventory.market.menu"
$ws.Range("C731").Value = "Tržiště"
$ws.Range("B732").Value = "inventory.root.home.menu"
$ws.Range("C732").Value = "Správa aplikace"
$ws.Range("B733").Value = "inventory.atomizer.index.title"
$ws.Range("C733").Value = "Vaše atomizéry"
$ws.Range("B734").Value = "inventory.liquid.menu"
$ws.Range("C734").Value = "Liquidy"
$ws.Range("B735").Value = "inventory.hardware.menu"
$ws.Range("C735").Value = "Hardware"

# Copy the formatting (style s="1", wrap text) from the last pre-existing row down onto the new rows
$ws.Range("A718:C718").Copy()
$ws.Range("A719:C735").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to the new last row, matching the authored view state
$ws.Range("B729").Select()

# Make "Translations - Common" the active sheet/tab (was "tabs" before)
$ws.Activate()

# Re-assert the selection on the now-active sheet (Activate can reset it)
$ws.Range("B729").Select()
